# Estado de Cuenta (NIT-9012814274) - refresh with new extract:
#  - Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Rebuilds the detail table (rows 16-39) with the new set of workers /
# periods, moves the signature footer down to rows 44-45, and refreshes
# the summary counters (VALOR MORA, Cant. Trabajadores, Cant. Periodos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ------------------------------------------------------------------
# 1. Preserve the two "special" row styles before we overwrite them:
#    - the bottom-bordered "last detail row" look currently on row 30
#    - the signature-footer look currently on rows 35-36
# ------------------------------------------------------------------
$ws.Range("B30:J30").Copy()
$ws.Range("B39:J39").PasteSpecial($xlPasteFormats)

$ws.Range("B35:C35").Copy()
$ws.Range("B44:C44").PasteSpecial($xlPasteFormats)
$ws.Range("H35:J35").Copy()
$ws.Range("H44:J44").PasteSpecial($xlPasteFormats)

$ws.Range("B36:C36").Copy()
$ws.Range("B45:C45").PasteSpecial($xlPasteFormats)
$ws.Range("H36:J36").Copy()
$ws.Range("H45:J45").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Row 30 is no longer the last detail row -> give it the normal
#    (non-bottom-bordered) look, then stamp that same look on the new
#    rows 31-38 that the refreshed extract needs.
# ------------------------------------------------------------------
$ws.Range("B29:J29").Copy()
$ws.Range("B30:J30").PasteSpecial($xlPasteFormats)
for ($r = 31; $r -le 38; $r++) {
  $ws.Range("B29:J29").Copy()
  $ws.Range("B$r`:J$r").PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Write the refreshed detail rows (16-39): Tipo Doc, No Doc,
#    Nombre, Periodo Mora, Valor Mora, Salario Basico.
# ------------------------------------------------------------------
$data = @(
  @("CC","9145360","JULIO ENRIQUE BELTRAN LLORENTE","2109",31125,1667436),
  @("CC","1047461453","MISAEL FARID TORRES MARTINEZ","2507",82420,2060500),
  @("CC","1001967527","ELVIRA YULIET MARTINEZ PEREZ","2104",10000,1500000),
  @("CC","1001967527","ELVIRA YULIET MARTINEZ PEREZ","2103",60000,1500000),
  @("CC","1001967527","ELVIRA YULIET MARTINEZ PEREZ","2102",60000,1500000),
  @("CC","1001967527","ELVIRA YULIET MARTINEZ PEREZ","2101",60000,1500000),
  @("CC","1001967527","ELVIRA YULIET MARTINEZ PEREZ","2012",60000,1500000),
  @("CC","1001967527","ELVIRA YULIET MARTINEZ PEREZ","2011",60000,1500000),
  @("CC","1001967527","ELVIRA YULIET MARTINEZ PEREZ","2010",60000,1500000),
  @("CC","1050950697","PEDRO LUIS LEON CARRILLO","2112",36341,908526),
  @("CC","1143380652","ANYIS JHOHANA MARTINEZ LOBO","2103",14000,1050000),
  @("CC","1143380652","ANYIS JHOHANA MARTINEZ LOBO","2102",42000,1050000),
  @("CC","1143380652","ANYIS JHOHANA MARTINEZ LOBO","2101",42000,1050000),
  @("CC","1143380652","ANYIS JHOHANA MARTINEZ LOBO","2012",42000,1050000),
  @("CC","1143380652","ANYIS JHOHANA MARTINEZ LOBO","2011",42000,1050000),
  @("CC","1143380652","ANYIS JHOHANA MARTINEZ LOBO","2010",42000,1050000),
  @("CC","1047512318","JUNIOR GERONIMO CABARCAS SOSCU","2009",35112,877803),
  @("PPT","5974130","LISBETH NAKARIT VERA RUIZ","2507",56940,1423500),
  @("CC","30655077","MARIA BERNARDA MORA PADILLA","2507",52000,1300000),
  @("CC","30655077","MARIA BERNARDA MORA PADILLA","2506",52000,1300000),
  @("CC","30655077","MARIA BERNARDA MORA PADILLA","2505",52000,1300000),
  @("CC","30655077","MARIA BERNARDA MORA PADILLA","2504",52000,1300000),
  @("CC","30655077","MARIA BERNARDA MORA PADILLA","2503",52000,1300000),
  @("CC","30655077","MARIA BERNARDA MORA PADILLA","2502",52000,1300000)
)

$r = 16
foreach ($row in $data) {
  $ws.Cells.Item($r, 2).Value = $row[0]
  $ws.Cells.Item($r, 3).Value = $row[1]
  $ws.Cells.Item($r, 4).Value = $row[2]
  $ws.Cells.Item($r, 5).Value = $row[3]
  $ws.Cells.Item($r, 6).Value = $row[4]
  $ws.Cells.Item($r, 7).Value = $row[5]
  $r = $r + 1
}

# ------------------------------------------------------------------
# 4. Move the signature footer text down to its new home (rows 44-45)
#    and clear it out of the old spot (now a plain detail row).
# ------------------------------------------------------------------
$ws.Cells.Item(44, 2).Value = "___________________________________"
$ws.Cells.Item(44, 8).Value = "___________________________________"
$ws.Cells.Item(45, 2).Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Cells.Item(45, 8).Value = "FIRMA DEL REPRESENTANTE LEGAL"

# ------------------------------------------------------------------
# 5. Refresh the summary counters.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 1147938
$ws.Range("C13").Value = 8
$ws.Range("F13").Value = 16

# ------------------------------------------------------------------
# 6. Merged cells for the relocated footer.
# ------------------------------------------------------------------
$ws.Range("B35:C35").UnMerge()
$ws.Range("H35:J35").UnMerge()
$ws.Range("B36:C36").UnMerge()
$ws.Range("H36:J36").UnMerge()
$ws.Range("B44:C44").Merge()
$ws.Range("H44:J44").Merge()
$ws.Range("B45:C45").Merge()
$ws.Range("H45:J45").Merge()
